$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '40.264.27'
$ws.Range("E2").Value = '  +0.25%  '

$ws.Range("D3").Value = '2.233.98'
$ws.Range("E3").Value = '  -0.46%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '295.14'
$ws.Range("E5").Value = '  +0.10%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '88.92'
$ws.Range("E6").Value = '  +2.38%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.513'
$ws.Range("E7").Value = '  -0.69%  '

$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.482'
$ws.Range("E9").Value = '  +1.70%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '30.47'
$ws.Range("E10").Value = '  -2.37%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0786'
$ws.Range("E11").Value = '  -1.70%  '

$ws.Range("E12").Value = '  +3.23%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.55'
$ws.Range("E13").Value = '  +1.18%  '

$ws.Range("D14").Value = '2.581.92'
$ws.Range("E14").Value = '  -0.43%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '13.90'
$ws.Range("E15").Value = '  -2.40%  '

$ws.Range("D16").Value = '2.220.25'
$ws.Range("E16").Value = '  -0.49%  '

$ws.Range("E17").Value = '  -0.19%  '

$ws.Range("D18").Value = '40.187.27'
$ws.Range("E18").Value = '  +0.23%  '

$ws.Range("E19").Value = '  -0.39%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.47'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.84'
$ws.Range("E21").Value = '  -0.07%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.71'
$ws.Range("E22").Value = '  -0.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.17'
$ws.Range("E23").Value = '  +0.18%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("E24").Value = '  -0.14%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.48'
$ws.Range("E25").Value = '  +0.30%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.82'
$ws.Range("E26").Value = '  -1.32%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.84'
$ws.Range("E27").Value = '  -0.82%  '

$ws.Range("E28").Value = '  -1.03%  '

$ws.Range("E29").Value = '  +0.16%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '155.12'
$ws.Range("E30").Value = '  +0.81%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.32'
$ws.Range("E31").Value = '  -3.22%  '

$ws.Range("E32").Value = '  -0.10%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.96'
$ws.Range("E33").Value = '  +1.09%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0721'
$ws.Range("E34").Value = '  +0.26%  '

$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.36'
$ws.Range("E35").Value = '  -1.36%  '

$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.91'
$ws.Range("E36").Value = '  +6.90%  '

$ws.Range("E37").Value = '  +0.46%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '15.87'
$ws.Range("E38").Value = '  -4.36%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0983'
$ws.Range("E39").Value = '  -2.40%  '

$ws.Range("E40").Value = '  -0.04%  '

$ws.Range("D41").Value = '2.130.34'
$ws.Range("E41").Value = '  +5.21%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.85'
$ws.Range("E42").Value = '  +0.90%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.14'
$ws.Range("E43").Value = '  -3.24%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '18.18'
$ws.Range("E44").Value = '  +11.28%  '

$ws.Range("E45").Value = '  -0.75%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.79'
$ws.Range("E46").Value = '  -1.78%  '

$ws.Range("E47").Value = '  +5.08%  '

$ws.Range("D48").Value = '2.448.40'
$ws.Range("E48").Value = '  -0.98%  '

$ws.Range("E49").Value = '  +1.62%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '69.48'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '89.14'
$ws.Range("E51").Value = '  -0.63%  '
